$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Title fix: "Ferbedienung" -> "Fernbedienung", and
#    " Steuern & Hinderniserkennung" -> " " + "s" + "teuern & Hinderniserkennung"
#    (also drops the now-obsolete spell-check proofErr markers around
#    the old misspelling)
# -----------------------------------------------------------------
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ferbedienung*") {
        $titlePara = $p
    }
}
# Replace the whole title paragraph's content (but not its paragraph
# mark) so the obsolete <w:proofErr/> markers around "Ferbedienung" are
# dropped along with the misspelling itself.
$titleRng = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)
$titleXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Titel"/></w:pPr><w:r><w:t xml:space="preserve">Arduino mit </w:t></w:r><w:r><w:t>Fernbedienung</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t>teuern &amp; Hinderniserkennung</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titleRng.InsertXML($titleXml) | Out-Null

# -----------------------------------------------------------------
# 2) Insert a new chapter heading "Das habe ich gelernt" right after
#    "Funktionsweise der Hardware", separated by an empty paragraph,
#    while leaving the pre-existing trailing empty heading paragraph
#    untouched.
# -----------------------------------------------------------------
$emptyXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Funktionsweise der Hardware`r") {
        $target = $p
    }
}

# Create a fresh blank paragraph right after "Funktionsweise der Hardware"
# and strip it down to a bare <w:p/>.
$target.Range.InsertParagraphAfter() | Out-Null
$blank = $d.Paragraphs($target.Index + 1)
$blankRng = $d.Range($blank.Range.Start, $blank.Range.End)
$blankRng.InsertXML($emptyXml) | Out-Null

# Create another fresh paragraph right after that blank one, and turn it
# into the new heading.
$blank2 = $d.Paragraphs($target.Index + 1)
$blank2.Range.InsertParagraphAfter() | Out-Null
$headingPara = $d.Paragraphs($blank2.Index + 1)
$headingRng = $d.Range($headingPara.Range.Start, $headingPara.Range.End)
$headingXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Das habe ich gelernt</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$headingRng.InsertXML($headingXml) | Out-Null
